# Applies the edits described by the commit "Changed subtitle line; Edited
# responsibility descriptions" to the active document.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2) | Out-Null
}

# --- Subtitle line (Heading2 under name) ---
Replace-Text "Language Scientist.Polymath.Data Expert.Cat Lover." `
             "Language Scientist. ML/Data Expert.Stress Baker."

# --- Villanova University (Assistant Professor title line) ---
# "Villanova U" is itself a uniformly-formatted (bold+italic) run, so
# extending it in place keeps that formatting intact.
Replace-Text "Villanova U" "Villanova University"

# --- Assistant Professor bullet list ---
Replace-Text "Develop and teach linguistics, cognitive science, AI, and Spanish courses " `
             "Design courses and instruct undergrads on linguistics, cognitive science, and AI"

Replace-Text "Advance an independent research program on individual differences in language use, cognitive processing, and human adaptation to change" `
             "Build an interdisciplinary, independent research program focused on language use, cognitive processing, and the human ability to adapt to change"

Replace-Text "Publish results of quantitative studies in" `
             "Publish results of quantitative research in"

Replace-Text "Mentor >20 undergraduate researchers in Language Use and Variation Lab" `
             "Mentor over 20 undergraduate researchers as principal investigator of LUV Lab"

Replace-Text "Write coherent, persuasive grant proposals to secure >`$50k in research funding" `
             "Write technical grant proposals to secure over `$50k in research funding"

# --- Applied Scientist II | Amazon (drop "Localization Tech") ---
Replace-Text "Applied Scientist II | Amazon Localization Tech" `
             "Applied Scientist II | Amazon"

Replace-Text "Audited data quality program for Prime Video subtitle localization " `
             "Fine-tuned a data quality program for Prime Video subtitle localization "

# --- Technical Program Manager II | Amazon (drop "Applied Modeling & Data Science") ---
Replace-Text "Technical Program Mgr. II | Amazon Applied Modeling & Data Science" `
             "Technical Program Manager II | Amazon"

Replace-Text "Owned a data quality program to support sourcing language data and building large language models for natural language processing" `
             "Created and managed a data quality program to support sourcing language data for internationalizing large language models using MTurk, SQL, and Tableau"

Replace-Text "Specified and monitored data quality alarms using AWS Cloudwatch" `
             "Specified KPIs and created automated alarms using AWS Cloudwatch"

Replace-Text "Forecasted data collection using generalized additive modeling" `
             "Forecasted data throughput using generalized additive modeling"

Replace-Text "Reported quality metrics to executives and stakeholders weekly" `
             "Compiled quality metrics and reports for executives and stakeholders weekly"

# --- Language Engineer II | Amazon (drop "Alexa Household Organization") ---
Replace-Text "Language Engineer II | Amazon Alexa Household Organization " `
             "Language Engineer II | Amazon"

Replace-Text "artifacts using git, FSTs, and virtual machines to production environments" `
             "artifacts to production environments using git, FSTs, and virtual machines"

# --- Selected Media Outreach heading ---
Replace-Text "Media Outreach" "Selected Media Outreach"

# --- Footer "Last Updated" date ---
$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.Find.Execute("19-Jun-2024", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "11-July-2024", 2) | Out-Null
